# issue #5: stock data output to json file
#
# Insert a new "property_category" column into the 股票 (stock) worksheet,
# right after the "total" column and before the "date" column. Every data
# row gets the literal value "stock" in that new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Column G = total, H = date (old layout). Insert a fresh column at H so
# everything from the old H onward (date, legislator_name, legislator_id)
# shifts one slot to the right.
$ws.Columns.Item(8).Insert()

# Header label for the freshly inserted column.
$ws.Range("H1").Value = "property_category"

# Find the last used row in column A (the numeric id column) so we fill
# every data row, however many there are.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
